# deactivation - november cycle
# Update the TestCases count and Instance name on the Config sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 (TestCases): 42 -> 40
$ws.Range("B2").Value = "40"

# D2 (Instance): Automation2 -> Automation1
$ws.Range("D2").Value = "Automation1"
